# Insert a new weekly price record at row 47 (pushing existing rows 47-96
# down to 48-97), matching the "Fruta / hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 47:96 down to 48:97, leaving a blank row at 47.
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with the new weekly observation.
$ws.Range("A47").Value = 6
$ws.Range("B47").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C47").Value = "Metropolitana"
$ws.Range("D47").Value = 44789
$ws.Range("E47").Value = 13
$ws.Range("F47").Value = 100114007
$ws.Range("G47").Value = "Jengibre"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 580
$ws.Range("K47").Value = 11000
$ws.Range("L47").Value = 12000
$ws.Range("M47").Value = 11448
$ws.Range("N47").Value = "`$/caja 13 kilos"
$ws.Range("O47").Value = "Perú"
$ws.Range("P47").Value = 881
$ws.Range("Q47").Value = 13
$ws.Range("R47").Value = "Hortaliza"
